$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 2 new rows before row 537, shifting existing rows 537:630 down to 539:632
$ws.Rows("537:538").Insert()

# Row 537 (new): Forelle, Primera
$ws.Range("A537").Value = 11
$ws.Range("B537").Value = "Vega Monumental Concepción"
$ws.Range("C537").Value = "Bíobío"
$ws.Range("D537").Value = Get-Date -Year 2023 -Month 3 -Day 30 -Hour 0 -Minute 0 -Second 0
$ws.Range("D537").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E537").Value = 8
$ws.Range("F537").Value = "Fruta"
$ws.Range("G537").Value = 100104
$ws.Range("H537").Value = "Frutos de pepita"
$ws.Range("I537").Value = 100104005
$ws.Range("J537").Value = "Pera"
$ws.Range("K537").Value = "Forelle"
$ws.Range("L537").Value = "Primera"
$ws.Range("M537").Value = 50
$ws.Range("N537").Value = 9000
$ws.Range("O537").Value = 9000
$ws.Range("P537").Value = 9000
$ws.Range("Q537").Value = "`$/caja 16 kilos empedrada"
$ws.Range("R537").Value = "Región de O'Higgins"
$ws.Range("S537").Value = 562
$ws.Range("T537").Value = 16

# Row 538 (new): Forelle, Segunda
$ws.Range("A538").Value = 11
$ws.Range("B538").Value = "Vega Monumental Concepción"
$ws.Range("C538").Value = "Bíobío"
$ws.Range("D538").Value = Get-Date -Year 2023 -Month 3 -Day 30 -Hour 0 -Minute 0 -Second 0
$ws.Range("D538").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E538").Value = 8
$ws.Range("F538").Value = "Fruta"
$ws.Range("G538").Value = 100104
$ws.Range("H538").Value = "Frutos de pepita"
$ws.Range("I538").Value = 100104005
$ws.Range("J538").Value = "Pera"
$ws.Range("K538").Value = "Forelle"
$ws.Range("L538").Value = "Segunda"
$ws.Range("M538").Value = 50
$ws.Range("N538").Value = 7000
$ws.Range("O538").Value = 7000
$ws.Range("P538").Value = 7000
$ws.Range("Q538").Value = "`$/caja 16 kilos empedrada"
$ws.Range("R538").Value = "Región de O'Higgins"
$ws.Range("S538").Value = 438
$ws.Range("T538").Value = 16
